$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.062216320625964
$ws.Range("C2").Value = 0.1141885135673419
$ws.Range("E2").Value = 0.1788703831046661
$ws.Range("F2").Value = 2.587394167838568
$ws.Range("G2").Value = 1.472043826745548
$ws.Range("H2").Value = 1.355903386789066
$ws.Range("J2").Value = 0.09545232200458287
$ws.Range("K2").Value = 0.547348162140878
$ws.Range("L2").Value = 0.3871886810256342
$ws.Range("N2").Value = 2.423160275878189
$ws.Range("B3").Value = 1.021508266728489
$ws.Range("C3").Value = 0.1130465409057102
$ws.Range("E3").Value = 0.1780122692324504
$ws.Range("F3").Value = 2.582276900787704
$ws.Range("G3").Value = 1.473784726153085
$ws.Range("H3").Value = 1.361684861232206
$ws.Range("J3").Value = 0.09557904080351065
$ws.Range("K3").Value = 0.5113466351630791
$ws.Range("L3").Value = 0.3797505931779597
$ws.Range("N3").Value = 2.444880360320575
$ws.Range("B4").Value = 0.9970280576508515
$ws.Range("C4").Value = 0.1123317331053784
$ws.Range("E4").Value = 0.1775485795231972
$ws.Range("F4").Value = 2.580332006271547
$ws.Range("G4").Value = 1.475605756212644
$ws.Range("H4").Value = 1.365757137647904
$ws.Range("J4").Value = 0.09566392449280592
$ws.Range("K4").Value = 0.4894860451655063
$ws.Range("L4").Value = 0.3753520217890696
$ws.Range("N4").Value = 2.458928697001607
$ws.Range("B5").Value = 0.9871821006285018
$ws.Range("C5").Value = 0.1120370066837495
$ws.Range("E5").Value = 0.1773755634505214
$ws.Range("F5").Value = 2.579840690564737
$ws.Range("G5").Value = 1.476536884493498
$ws.Range("H5").Value = 1.367548067151176
$ws.Range("J5").Value = 0.09570030299567023
$ws.Range("K5").Value = 0.4806394438443817
$ws.Range("L5").Value = 0.3736020763375194
$ws.Range("N5").Value = 2.464832549568172
$ws.Range("B6").Value = 0.9855550490536018
$ws.Range("C6").Value = 0.1119878597849429
$ws.Range("E6").Value = 0.1773477985533845
$ws.Range("F6").Value = 2.579777308842452
$ws.Range("G6").Value = 1.476702913097355
$ws.Range("H6").Value = 1.367853390748877
$ws.Range("J6").Value = 0.09570645182499504
$ws.Range("K6").Value = 0.4791742114866793
$ws.Range("L6").Value = 0.3733140710757112
$ws.Range("N6").Value = 2.465823692724129
$ws.Range("B7").Value = 0.996894744897304
$ws.Range("C7").Value = 0.1123277722461147
$ws.Range("E7").Value = 0.1775461815554777
$ws.Range("F7").Value = 2.580324160202991
$ws.Range("G7").Value = 1.475617548412984
$ws.Range("H7").Value = 1.365780758441844
$ws.Range("J7").Value = 0.09566440785722818
$ws.Range("K7").Value = 0.4893664863347738
$ws.Range("L7").Value = 0.3753282490970236
$ws.Range("N7").Value = 2.459007593530281
$ws.Range("B8").Value = 1.048073642684926
$ws.Range("C8").Value = 0.1137975852765081
$ws.Range("E8").Value = 0.1785614186529259
$ws.Range("F8").Value = 2.585381392122912
$ws.Range("G8").Value = 1.472487936103064
$ws.Range("H8").Value = 1.357788452947958
$ws.Range("J8").Value = 0.09549455038499133
$ws.Range("K8").Value = 0.5348842506236622
$ws.Range("L8").Value = 0.384589155213348
$ws.Range("N8").Value = 2.430501547106175
$ws.Range("B9").Value = 1.152503491077852
$ws.Range("C9").Value = 0.1165722123561608
$ws.Range("E9").Value = 0.181051773110859
$ws.Range("F9").Value = 2.604791359081958
$ws.Range("G9").Value = 1.472323722681736
$ws.Range("H9").Value = 1.346258113698255
$ws.Range("J9").Value = 0.09521727301672378
$ws.Range("K9").Value = 0.62607686344748
$ws.Range("L9").Value = 0.404080958731555
$ws.Range("N9").Value = 2.380252001217038
$ws.Range("B10").Value = 1.231696325665894
$ws.Range("C10").Value = 0.1185459231859127
$ws.Range("E10").Value = 0.1831836810863621
$ws.Range("F10").Value = 2.624836284913755
$ws.Range("G10").Value = 1.475854042320208
$ws.Range("H10").Value = 1.340309494104687
$ws.Range("J10").Value = 0.09504711812803635
$ws.Range("K10").Value = 0.6942515110039835
$ws.Range("L10").Value = 0.4192078805967157
$ws.Range("N10").Value = 2.346777907948216
$ws.Range("B11").Value = 1.268257238249646
$ws.Range("C11").Value = 0.1194299309011839
$ws.Range("E11").Value = 0.1842187188886157
$ws.Range("F11").Value = 2.635211112286612
$ws.Range("G11").Value = 1.478255121502613
$ws.Range("H11").Value = 1.338150625430913
$ws.Range("J11").Value = 0.09497690018578808
$ws.Range("K11").Value = 0.7255210831171723
$ws.Range("L11").Value = 0.4262636017144672
$ws.Range("N11").Value = 2.332297457356383
$ws.Range("B12").Value = 1.282178572249052
$ws.Range("C12").Value = 0.1197627025645147
$ws.Range("E12").Value = 0.1846199946348079
$ws.Range("F12").Value = 2.639320328769614
$ws.Range("G12").Value = 1.479278830912534
$ws.Range("H12").Value = 1.337411755807523
$ws.Range("J12").Value = 0.09495133597726912
$ws.Range("K12").Value = 0.7373987769698829
$ws.Range("L12").Value = 0.42896037525243
$ws.Range("N12").Value = 2.326921558830975
$ws.Range("B13").Value = 1.279176969394143
$ws.Range("C13").Value = 0.1196911222973682
$ws.Range("E13").Value = 0.184533158582596
$ws.Range("F13").Value = 2.638427310444357
$ws.Range("G13").Value = 1.479053263391052
$ws.Range("H13").Value = 1.337567387338879
$ws.Range("J13").Value = 0.09495679618435204
$ws.Range("K13").Value = 0.7348390808685963
$ws.Range("L13").Value = 0.4283784707305216
$ws.Range("N13").Value = 2.328074570467891
$ws.Range("B14").Value = 1.269401024391982
$ws.Range("C14").Value = 0.1194573479215535
$ws.Range("E14").Value = 0.1842515453938489
$ws.Range("F14").Value = 2.635545563524488
$ws.Range("G14").Value = 1.478337047705338
$ws.Range("H14").Value = 1.338088262189842
$ws.Range("J14").Value = 0.0949747764870299
$ws.Range("K14").Value = 0.7264975366878161
$ws.Range("L14").Value = 0.4264849682458305
$ws.Range("N14").Value = 2.331853023200672
$ws.Range("B15").Value = 1.263422924166036
$ws.Range("C15").Value = 0.1193138963694977
$ws.Range("E15").Value = 0.1840802626802862
$ws.Range("F15").Value = 2.633803911358484
$ws.Range("G15").Value = 1.477913256704582
$ws.Range("H15").Value = 1.338417554147071
$ws.Range("J15").Value = 0.09498592330782074
$ws.Range("K15").Value = 0.7213928553244671
$ws.Range("L15").Value = 0.4253283846945806
$ws.Range("N15").Value = 2.33418144183841
$ws.Range("B16").Value = 1.229317718911261
$ws.Range("C16").Value = 0.1184878733039838
$ws.Range("E16").Value = 0.1831173471368039
$ws.Range("F16").Value = 2.624183531740712
$ws.Range("G16").Value = 1.475713137549405
$ws.Range("H16").Value = 1.340461588803265
$ws.Range("J16").Value = 0.0950518509936753
$ws.Range("K16").Value = 0.6922131074495326
$ws.Range("L16").Value = 0.4187502689485001
$ws.Range("N16").Value = 2.347739281489723
$ws.Range("B17").Value = 1.20853211237511
$ws.Range("C17").Value = 0.1179775948323964
$ws.Range("E17").Value = 0.182543299364081
$ws.Range("F17").Value = 2.618603393358342
$ws.Range("G17").Value = 1.474567180963589
$ws.Range("H17").Value = 1.341855658481066
$ws.Range("J17").Value = 0.09509413067843653
$ws.Range("K17").Value = 0.6743777729597298
$ws.Range("L17").Value = 0.4147593681017128
$ws.Range("N17").Value = 2.356247995767468
$ws.Range("B18").Value = 1.196627234540756
$ws.Range("C18").Value = 0.117682792692591
$ws.Range("E18").Value = 0.1822192642395635
$ws.Range("F18").Value = 2.615512102703434
$ws.Range("G18").Value = 1.473982889219812
$ws.Range("H18").Value = 1.342708996597935
$ws.Range("J18").Value = 0.09511912570011027
$ws.Range("K18").Value = 0.6641435296273244
$ws.Range("L18").Value = 0.412480332877962
$ws.Range("N18").Value = 2.361212295607316
$ws.Range("B19").Value = 1.192605128800096
$ws.Range("C19").Value = 0.1175827536470706
$ws.Range("E19").Value = 0.1821106079126444
$ws.Range("F19").Value = 2.614485760471609
$ws.Range("G19").Value = 1.473797906463062
$ws.Range("H19").Value = 1.343006769628602
$ws.Range("J19").Value = 0.09512770506926849
$ws.Range("K19").Value = 0.660682551272771
$ws.Range("L19").Value = 0.4117115165016401
$ws.Range("N19").Value = 2.362905195352589
$ws.Range("B20").Value = 1.210739557736701
$ws.Range("C20").Value = 0.1180320496930847
$ws.Range("E20").Value = 0.1826037724176572
$ws.Range("F20").Value = 2.619185170570731
$ws.Range("G20").Value = 1.474681424166349
$ws.Range("H20").Value = 1.341701926970032
$ws.Range("J20").Value = 0.0950895599358681
$ws.Range("K20").Value = 0.6762738755816997
$ws.Range("L20").Value = 0.415182507465687
$ws.Range("N20").Value = 2.355334950782563
$ws.Range("B21").Value = 1.272270384112687
$ws.Range("C21").Value = 0.1195260668750819
$ws.Range("E21").Value = 0.1843340092105734
$ws.Range("F21").Value = 2.636387105136833
$ws.Range("G21").Value = 1.478544309766434
$ws.Range("H21").Value = 1.337933134455568
$ws.Range("J21").Value = 0.09496946745693791
$ws.Range("K21").Value = 0.7289466599682726
$ws.Range("L21").Value = 0.4270404605965723
$ws.Range("N21").Value = 2.330740280754696
$ws.Range("B22").Value = 1.31293005274199
$ws.Range("C22").Value = 0.1204909394096632
$ws.Range("E22").Value = 0.185519174941529
$ws.Range("F22").Value = 2.648681551296875
$ws.Range("G22").Value = 1.481736227870925
$ws.Range("H22").Value = 1.33592840950142
$ws.Range("J22").Value = 0.09489695575672386
$ws.Range("K22").Value = 0.7635844742450786
$ws.Range("L22").Value = 0.4349355292704615
$ws.Range("N22").Value = 2.315293082656908
$ws.Range("B23").Value = 1.291188611551206
$ws.Range("C23").Value = 0.119977023394739
$ws.Range("E23").Value = 0.1848816724292064
$ws.Range("F23").Value = 2.642023569124831
$ws.Range("G23").Value = 1.479971538076896
$ws.Range("H23").Value = 1.336956437587418
$ws.Range("J23").Value = 0.09493511233486096
$ws.Range("K23").Value = 0.7450782242803484
$ws.Range("L23").Value = 0.4307085467778222
$ws.Range("N23").Value = 2.323480151650934
$ws.Range("B24").Value = 1.209741431340376
$ws.Range("C24").Value = 0.1180074351252998
$ws.Range("E24").Value = 0.1825764138807777
$ws.Range("F24").Value = 2.618921785235955
$ws.Range("G24").Value = 1.474629542650518
$ws.Range("H24").Value = 1.341771267408674
$ws.Range("J24").Value = 0.09509162422502371
$ws.Range("K24").Value = 0.6754165867358779
$ws.Range("L24").Value = 0.4149911582209
$ws.Range("N24").Value = 2.355747512500088
$ws.Range("B25").Value = 1.123818153692184
$ws.Range("C25").Value = 0.1158330434299231
$ws.Range("E25").Value = 0.1803248268150384
$ws.Range("F25").Value = 2.598524460054932
$ws.Range("G25").Value = 1.471727704208291
$ws.Range("H25").Value = 1.348934171101163
$ws.Range("J25").Value = 0.09528635694475796
$ws.Range("K25").Value = 0.6012002411365813
$ws.Range("L25").Value = 0.3986659382875217
$ws.Range("N25").Value = 2.393240960532616
